# dedication.xlsx: split the single "Hours dedicated" column (B) into three
# identical columns (B, C, D) -- mirroring the HTML-divided-into-React-components
# change, each column now rendered by its own (duplicated) component -- and add a
# new shared string "asdf" used as a (style-less) header label for the two new
# columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Give the existing "Hours dedicated" values (B2:B7) an underlined font.
#    This is also what introduces the new font (fontId 2) and cell style
#    (cellXfs index 3) used later by the duplicated columns.
$ws.Range("B2:B7").Font.Underline = $true

# 2. New header cells for the two duplicated columns (plain/no style, like the
#    rest of row 1 besides A1/B1) -- both reuse the same new shared string "asdf".
$ws.Range("C1").Value = "asdf"
$ws.Range("D1").Value = "asdf"

# 3. Duplicate the numeric values from column B into the new columns C and D.
for ($r = 2; $r -le 7; $r++) {
    $val = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 3).Value = $val
    $ws.Cells.Item($r, 4).Value = $val
}

# 4. Copy the (now underlined, centered) formatting from B2:B7 onto C2:D7 so
#    they share the exact same cell style, without leaving any stray/unused
#    style entries behind.
$ws.Range("B2:B7").Copy()
$ws.Range("C2:D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 5. Leave the sheet with the same selection as the saved workbook.
$null = $ws.Range("G11").Select()
